$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.137.43'
Set-TextCell $ws 'E2' '  -3.47%  '
Set-TextCell $ws 'D3' '1.849.75'
Set-TextCell $ws 'E3' '  -2.42%  '
Set-TextCell $ws 'E4' '  -0.30%  '
Set-TextCell $ws 'D5' '0.7063'
Set-TextCell $ws 'E5' '  -4.69%  '
Set-TextCell $ws 'D6' '238.28'
Set-TextCell $ws 'E6' '  -2.04%  '
Set-TextCell $ws 'E7' '  -0.25%  '
Set-TextCell $ws 'D8' '0.3041'
Set-TextCell $ws 'E8' '  -4.25%  '
Set-TextCell $ws 'D9' '0.07478'
Set-TextCell $ws 'E9' '  +3.43%  '
Set-TextCell $ws 'D10' '23.38'
Set-TextCell $ws 'D11' '0.08132'
Set-TextCell $ws 'E11' '  -2.75%  '
Set-TextCell $ws 'B12' 'WrappedEther'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D12' '1.860.11'
Set-TextCell $ws 'E12' '  -4.25%  '
Set-TextCell $ws 'B13' 'Polygon'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws 'D13' '0.7257'
Set-TextCell $ws 'E13' '  -4.80%  '
Set-TextCell $ws 'D14' '5.211'
Set-TextCell $ws 'E14' '  -4.63%  '
Set-TextCell $ws 'D15' '88.88'
Set-TextCell $ws 'E15' '  -4.56%  '
Set-TextCell $ws 'D16' '29.188.38'
Set-TextCell $ws 'E16' '  -3.48%  '
Set-TextCell $ws 'D17' '5.759'
Set-TextCell $ws 'E17' '  -6.82%  '
Set-TextCell $ws 'D18' '238.80'
Set-TextCell $ws 'E18' '  -4.91%  '
Set-TextCell $ws 'D19' '13.08'
Set-TextCell $ws 'E19' '  -4.28%  '
Set-TextCell $ws 'D20' '0.000007652'
Set-TextCell $ws 'E20' '  -3.19%  '
Set-TextCell $ws 'E21' '  -0.32%  '
Set-TextCell $ws 'D22' '2.102.97'
Set-TextCell $ws 'E22' '  -4.54%  '
Set-TextCell $ws 'D23' '1.000'
Set-TextCell $ws 'E23' '  -0.29%  '
Set-TextCell $ws 'D24' '7.546'
Set-TextCell $ws 'E24' '  -5.51%  '
Set-TextCell $ws 'D25' '161.90'
Set-TextCell $ws 'E25' '  -1.84%  '
Set-TextCell $ws 'D26' '8.990'
Set-TextCell $ws 'E26' '  -3.58%  '
Set-TextCell $ws 'D27' '0.1459'
Set-TextCell $ws 'E27' '  -8.02%  '
Set-TextCell $ws 'D28' '18.03'
Set-TextCell $ws 'E28' '  -4.15%  '
Set-TextCell $ws 'D29' '1.951'
Set-TextCell $ws 'E29' '  -5.76%  '
Set-TextCell $ws 'D30' '1.389'
Set-TextCell $ws 'E30' '  -6.63%  '
Set-TextCell $ws 'D31' '4.510'
Set-TextCell $ws 'E31' '  -1.92%  '
Set-TextCell $ws 'D32' '1.493'
Set-TextCell $ws 'E32' '  -2.79%  '
Set-TextCell $ws 'D33' '3.980'
Set-TextCell $ws 'E33' '  -5.55%  '
Set-TextCell $ws 'D34' '0.05153'
Set-TextCell $ws 'E34' '  -4.36%  '
Set-TextCell $ws 'E35' '  -5.70%  '
Set-TextCell $ws 'D36' '1.034'
Set-TextCell $ws 'E36' '  +3.05%  '
Set-TextCell $ws 'D37' '0.6995'
Set-TextCell $ws 'E37' '  -11.49%  '
Set-TextCell $ws 'D38' '2.654'
Set-TextCell $ws 'E38' '  -3.00%  '
Set-TextCell $ws 'D39' '0.01876'
Set-TextCell $ws 'E39' '  -4.89%  '
Set-TextCell $ws 'D40' '2.676'
Set-TextCell $ws 'E40' '  -3.32%  '
Set-TextCell $ws 'D41' '0.9342'
Set-TextCell $ws 'E41' '  +6.83%  '
Set-TextCell $ws 'D42' '1.084.26'
Set-TextCell $ws 'E42' '  -1.56%  '
Set-TextCell $ws 'D43' '5.996'
Set-TextCell $ws 'E43' '  -1.58%  '
Set-TextCell $ws 'D44' '0.4287'
Set-TextCell $ws 'E44' '  -6.38%  '
Set-TextCell $ws 'D45' '69.89'
Set-TextCell $ws 'E45' '  -4.27%  '
Set-TextCell $ws 'E46' '  -0.40%  '
Set-TextCell $ws 'D47' '102.23'
Set-TextCell $ws 'E47' '  -2.31%  '
Set-TextCell $ws 'E48' '  -6.84%  '
Set-TextCell $ws 'D49' '2.000.49'
Set-TextCell $ws 'E49' '  -3.73%  '
Set-TextCell $ws 'D50' '9.165'
Set-TextCell $ws 'E50' '  -5.06%  '
Set-TextCell $ws 'E51' '  -7.64%  '
